$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 39.447365
$ws.Range("H2").Value = 118.342095
$ws.Range("I2").Value = 0.381567221408747
$ws.Range("J2").Value = 0.3841102728938004
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.642212666666667
$ws.Range("N2").Value = 4.926638000000001
$ws.Range("Q2").Value = 64.78096246962333
$ws.Range("R2").Value = 583.0286622266101
$ws.Range("S2").Value = 0.381567221408747
$ws.Range("T2").Value = 0.3841102728938004

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 54.64926166666667
$ws.Range("H3").Value = 163.947785
$ws.Range("I3").Value = 0.5286124162206919
$ws.Range("J3").Value = 0.5321354876866436
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.642212666666667
$ws.Range("N3").Value = 4.926638000000001
$ws.Range("Q3").Value = 89.74570973298111
$ws.Range("R3").Value = 807.7113875968301
$ws.Range("S3").Value = 0.5286124162206919
$ws.Range("T3").Value = 0.5321354876866436

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.463205333333333
$ws.Range("H4").Value = 10.389616
$ws.Range("I4").Value = 0.03349895832606192
$ws.Range("J4").Value = 0.03372222062674989
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.642212666666667
$ws.Range("N4").Value = 4.926638000000001
$ws.Range("Q4").Value = 5.687319665667555
$ws.Range("R4").Value = 51.18587699100801
$ws.Range("S4").Value = 0.03349895832606192
$ws.Range("T4").Value = 0.03372222062674989

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 3.769275333333333
$ws.Range("H5").Value = 11.307826
$ws.Range("I5").Value = 0.03645951803534987
$ws.Range("J5").Value = 0.03670251173680516
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.642212666666667
$ws.Range("N5").Value = 4.926638000000001
$ws.Range("Q5").Value = 6.189951696554223
$ws.Range("R5").Value = 55.70956526898801
$ws.Range("S5").Value = 0.03645951803534987
$ws.Range("T5").Value = 0.03670251173680516

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("G6").Value = 2.053371
$ws.Range("H6").Value = 4.106742
$ws.Range("I6").Value = 0.01986188600914916
$ws.Range("J6").Value = 0.01332950705600092
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.642212666666667
$ws.Range("N6").Value = 4.926638000000001
$ws.Range("Q6").Value = 3.372071865566
$ws.Range("R6").Value = 20.232431193396
$ws.Range("S6").Value = 0.01986188600914916
$ws.Range("T6").Value = 0.01332950705600092
